# Natmi following Dr Hou advice
# Re-run of the LR-pair (Tnfsf11 -> Tnfrsf11b) analysis across all
# sender/receiver cluster combinations (FAPs, M2, sCs), replacing the
# single-pair row 2 with the full 3x3 combination (rows 2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$header = @(
    "Sending cluster","Ligand symbol","Receptor symbol","Target cluster",
    "Ligand-expressing cells","Ligand detection rate",
    "Ligand average expression value","Ligand total expression value",
    "Ligand derived specificity of average expression value",
    "Ligand derived specificity of total expression value",
    "Receptor-expressing cells","Receptor detection rate",
    "Receptor average expression value","Receptor total expression value",
    "Receptor derived specificity of average expression value",
    "Receptor derived specificity of total expression value",
    "Edge average expression weight","Edge total expression weight",
    "Edge average expression derived specificity",
    "Edge total expression derived specificity"
)

# One row per (sending cluster, target cluster) combination for the
# Tnfsf11-Tnfrsf11b ligand-receptor pair.
$rows = @(
    @("FAPs","Tnfsf11","Tnfrsf11b","FAPs",3,1,1.690534333333333,5.071603,0.9339531815638826,0.9339531815638826,3,1,3.776574666666666,11.329724,0.9855052394405499,0.9855052394405499,6.384429136396887,57.45986222757199,0.9204157538233775,0.9204157538233775),
    @("FAPs","Tnfsf11","Tnfrsf11b","sCs",3,1,1.690534333333333,5.071603,0.9339531815638826,0.9339531815638826,1,0.3333333333333333,0.05554566666666667,0.166637,0.01449476055945007,0.01449476055945008,0.09390185656788888,0.845116709111,0.01353742774050508,0.01353742774050508),
    @("M2","Tnfsf11","Tnfrsf11b","FAPs",1,0.3333333333333333,0.014756,0.044268,0.008152104855500313,0.008152104855500315,3,1,3.776574666666666,11.329724,0.9855052394405499,0.9855052394405499,0.05572713578133333,0.501544222032,0.008033942047564305,0.008033942047564307),
    @("M2","Tnfsf11","Tnfrsf11b","sCs",1,0.3333333333333333,0.014756,0.044268,0.008152104855500313,0.008152104855500315,1,0.3333333333333333,0.05554566666666667,0.166637,0.01449476055945007,0.01449476055945008,0.0008196318573333334,0.007376686716000001,0.0001181628079360074,0.0001181628079360074),
    @("sCs","Tnfsf11","Tnfrsf11b","FAPs",1,0.3333333333333333,0.1047943333333333,0.314383,0.05789471358061704,0.05789471358061705,3,1,3.776574666666666,11.329724,0.9855052394405499,0.9855052394405499,0.3957636244768889,3.561872620292,0.05705554356960806,0.05705554356960806),
    @("sCs","Tnfsf11","Tnfrsf11b","sCs",1,0.3333333333333333,0.1047943333333333,0.314383,0.05789471358061704,0.05789471358061705,1,0.3333333333333333,0.05554566666666667,0.166637,0.01449476055945007,0.01449476055945008,0.005820871107888889,0.05238783997100001,0.0008391700110089865,0.0008391700110089867)
)

# Header row (unchanged, re-asserted defensively).
for ($c = 1; $c -le $header.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $header[$c - 1]
}

# Data rows starting at row 2.
$r = 2
foreach ($row in $rows) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r = $r + 1
}
